# Apply weekly update: insert a new latest-week record at row 20 and push
# all subsequent historical records down by one row, appending the former
# last row (68) as new row 69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the last existing row (68) into the new row (69), copying every column.
#    Column D holds dates; read it with Value2 (raw serial number) and re-apply the
#    date number format explicitly so the brand-new cell doesn't pick up a default
#    style instead of reusing the sheet's existing date style.
for ($c = 1; $c -le 20; $c++) {
    if ($c -eq 4) {
        $v = $ws.Cells.Item(68, $c).Value2()
        $ws.Cells.Item(69, $c).Value = $v
        $ws.Cells.Item(69, $c).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    } else {
        $v = $ws.Cells.Item(68, $c).Value()
        $ws.Cells.Item(69, $c).Value = $v
    }
}

# 2) Shift rows 68 down to 21 into the row below (only the columns that vary
#    per record: D, L, M, N, O, P, Q, R, S, T). Columns A, B, C, E-K are
#    constant across all rows of this sheet, so they don't need shifting.
#    Process from the bottom up so we never overwrite a source row before
#    it has been read.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)
for ($r = 68; $r -ge 21; $r--) {
    foreach ($c in $cols) {
        $v = $ws.Cells.Item($r - 1, $c).Value()
        $ws.Cells.Item($r, $c).Value = $v
    }
}

# 3) Write the brand-new record into row 20.
$ws.Cells.Item(20, 4).Value = 45274
$ws.Cells.Item(20, 13).Value = 50
$ws.Cells.Item(20, 14).Value = 16000
$ws.Cells.Item(20, 15).Value = 16000
$ws.Cells.Item(20, 16).Value = 16000
$ws.Cells.Item(20, 19).Value = 1600

Write-Output "done"
